# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9
$ws1.Range("F3").Value = 151
$ws1.Range("F4").Value = 711

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9
$ws4.Range("F4").Value = 151
$ws4.Range("F5").Value = 711
